# Senior Project Tracking RTC - "Lots of updates. At 88.58 hours." commit.
# Adds 14 new log rows (rows 9-22), extends the D-column elapsed-hours
# formula down to row 22, and normalizes / restyles the Date column.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1. New data rows 9-22 (Date / Start time / End time / Reason).
#    Column D (elapsed hours) is filled in further down with one shared
#    formula block; L1 (=SUM(D:D)) recalculates automatically.
# ---------------------------------------------------------------------
$ws.Range("A9").Value = 42415
$ws.Range("B9").Value = 0.38541666666666669
$ws.Range("C9").Value = 0.62847222222222221
$ws.Range("E9").Value = "Adding component to ensure that dependencies are loaded correctly for RTC."

$ws.Range("A10").Value = 42416
$ws.Range("B10").Value = 0.375
$ws.Range("C10").Value = 0.58333333333333337
$ws.Range("E10").Value = "Various RTC related things"

$ws.Range("A11").Value = 42416
$ws.Range("B11").Value = 0.66666666666666663
$ws.Range("C11").Value = 0.70833333333333337
$ws.Range("E11").Value = "Meeting for RTC update and path forward for what we want to do and how we want to do it."

$ws.Range("A12").Value = 42417
$ws.Range("B12").Value = 0.375
$ws.Range("C12").Value = 0.70833333333333337
$ws.Range("E12").Value = "Page Navigation"

$ws.Range("A13").Value = 42418
$ws.Range("B13").Value = 0.375
$ws.Range("C13").Value = 0.70833333333333337
$ws.Range("E13").Value = "Page Navigation"

$ws.Range("A14").Value = 42419
$ws.Range("B14").Value = 0.375
$ws.Range("C14").Value = 0.70833333333333337
$ws.Range("E14").Value = "Page Navigation"

$ws.Range("A15").Value = 42422
$ws.Range("B15").Value = 0.375
$ws.Range("C15").Value = 0.5
$ws.Range("E15").Value = "Finished Page Navigation"

$ws.Range("A16").Value = 42422
$ws.Range("B16").Value = 0.5
$ws.Range("C16").Value = 0.625
$ws.Range("E16").Value = "Working on presentation 1"

$ws.Range("A17").Value = 42422
$ws.Range("B17").Value = 0.64583333333333337
$ws.Range("C17").Value = 0.6875
$ws.Range("E17").Value = "Presenting"

$ws.Range("A18").Value = 42423
$ws.Range("B18").Value = 0.375
$ws.Range("C18").Value = 0.66666666666666663
$ws.Range("E18").Value = "Working on stream pull functionality"

$ws.Range("A19").Value = 42424
$ws.Range("B19").Value = 0.375
$ws.Range("C19").Value = 0.70833333333333337
$ws.Range("E19").Value = "Working on stream pull functionality"

$ws.Range("A20").Value = 42425
$ws.Range("B20").Value = 0.375
$ws.Range("C20").Value = 0.70833333333333337
$ws.Range("E20").Value = "Working on stream pull functionality"

$ws.Range("A21").Value = 42426
$ws.Range("B21").Value = 0.28125
$ws.Range("C21").Value = 0.33333333333333331
$ws.Range("E21").Value = "Presenting initial demo of stream pull functionality"

$ws.Range("A22").Value = 42426
$ws.Range("B22").Value = 0.33333333333333331
$ws.Range("C22").Value = 0.375
$ws.Range("E22").Value = "Touch-ups and troubleshooting."

# ---------------------------------------------------------------------
# 2. Elapsed-hours formula (column D) for the new rows.
# ---------------------------------------------------------------------
$ws.Range("D9:D22").Formula = "=MOD(IF(ISBLANK(C9),B9, C9)-B9, 1)*24"

# ---------------------------------------------------------------------
# 3. Number formats for the new rows (matches the existing columns).
# ---------------------------------------------------------------------
$ws.Range("A9:A22").NumberFormat = "d-mmm"
$ws.Range("B9:C22").NumberFormat = "h:mm AM/PM"
$ws.Range("D9:D22").NumberFormat = "0.00"

# ---------------------------------------------------------------------
# 4. Normalize the whole Date column (A) to the standard "d-mmm" date
#    format, including row 8 which previously used a one-off format.
# ---------------------------------------------------------------------
$ws.Columns("A").NumberFormat = "d-mmm"

# ---------------------------------------------------------------------
# 5. Re-style the "Date" header (A1) to bold + underline, matching the
#    other section headers.
# ---------------------------------------------------------------------
$ws.Range("A1").Font.Bold = $true
$ws.Range("A1").Font.Underline = $true

# ---------------------------------------------------------------------
# 6. Tidy up the active selection.
# ---------------------------------------------------------------------
[void]$ws.Range("A1").Select()
